$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the same formatting used by the other headers (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Populate the new "Save" column values - 0 for every row except row 6 (=1)
$saveValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 1
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
